# Update the build timestamp embedded in the version string from
# "17.29.55 EST" to "18.05.36 EST" throughout the workbook.

$wb = $excel.ActiveWorkbook

$oldTime = "17.29.55 EST"
$newTime = "18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---------------------------------------------------
# A2: "Version: ... (built on February 03 2026 17.29.55 EST)"
$cell = $aboutSheet.Range("A2")
$v = $cell.Value()
$cell.Value = $v.Replace($oldTime, $newTime)

# A6: "Recommended Citation: ... (built on February 03 2026 17.29.55 EST)' ..."
$cell = $aboutSheet.Range("A6")
$v = $cell.Value()
$cell.Value = $v.Replace($oldTime, $newTime)

# --- "Boundaries and methane sources" sheet --------------------------
# Column S ("build_version") rows 2 through 33 each contain
# "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"
for ($row = 2; $row -le 33; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    $v = $cell.Value()
    if ($v -ne $null -and $v.ToString().Contains($oldTime)) {
        $cell.Value = $v.ToString().Replace($oldTime, $newTime)
    }
}
